# Regenerate save_data to use K instead of Strike# (column G), regen std/mean,
# calc and write s_vals. Only column G values change (rows 2-40).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2 = 2
    3 = 6
    4 = 5
    5 = 9
    6 = 4
    7 = 4
    8 = 6
    9 = 6
    10 = 0
    11 = 3
    12 = 1
    13 = 8
    14 = 4
    15 = 5
    16 = 6
    17 = 8
    18 = 5
    19 = 1
    20 = 7
    21 = 8
    22 = 6
    23 = 1
    24 = 4
    25 = 3
    26 = 4
    27 = 4
    28 = 3
    29 = 2
    30 = 3
    31 = 4
    32 = 4
    33 = 4
    34 = 4
    35 = 4
    36 = 6
    37 = 3
    38 = 2
    39 = 6
    40 = 1
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}
